# Add a new worksheet "Extra" after "Must-Have" containing the
# "READ - Users should be able to retrieve all contacts in an address book"
# screen operation (REST endpoint test cases), and tweak the view/selection
# state on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Tweak "Must-Have" sheet view: zoom 120 -> 140, selection C2 -> A2
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.Zoom = 140

# ---------------------------------------------------------------------
# 2. Create the new "Extra" worksheet right after "Must-Have"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Extra"

# Column widths matching the "Must-Have" look (col A wide, rest default)
$ws2.Columns.Item(1).ColumnWidth = 17.5

# ---------------------------------------------------------------------
# 3. Fill in the data - order chosen so that shared-string allocation
#    matches the target layout.
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Requirement"
$ws2.Range("B1").Value = "Case"
$ws2.Range("C1").Value = "Expected result"
$ws2.Range("D1").Value = "Actual result (dev)"
$ws2.Range("E1").Value = "Actual result (prod)"
$ws2.Range("F1").Value = "Remark"

$ws2.Range("A2").Value = "READ - Users should be able to retrieve all contacts in an address book`tRead the entries from all phone books"
$ws2.Range("C2").Value = "All entries from all phone books should be shown"
$ws2.Range("B2").Value = "Read the entries from all phone books"
$ws2.Range("D2").Value = "Passed"

$ws2.Range("A3").Value = "RESTful API endpoint for createEntry"
$ws2.Range("A4").Value = "RESTful API endpoint for updateEntry"
$ws2.Range("A5").Value = "RESTful API endpoint for readAllEntriesFromAllPhoneBooks"
$ws2.Range("A6").Value = "RESTful API endpoint for readAllEntriesFromSinglePhoneBook"
$ws2.Range("A7").Value = "RESTful API endpoint for readUniqueEntriesFromAllPhoneBooks"
$ws2.Range("A8").Value = "RESTful API endpoint for deleteEntryFromSinglePhoneBook"

$ws2.Range("C3").Value = "A message should be returned indicating the successful creation of the entry"
$ws2.Range("C4").Value = "A message should be returned indicating the successful update of the entry"
$ws2.Range("C5").Value = "A list containing all phone book entries should be returned"

$ws2.Range("B3").Value = "Create a new customer entry in the specified phone book"
$ws2.Range("B4").Value = "Update an existing customer entry in the specified phone book"

$ws2.Range("B6").Value = "Read the entries from the specified phone book "
$ws2.Range("C6").Value = "The phoneBook entries of the specified phone book should be returned"

$ws2.Range("B7").Value = "Read unique entries across all phone books"
$ws2.Range("C7").Value = "Unique entires from across all phone books should be returned"

$ws2.Range("B8").Value = "Delete an existing entry from the specified phone book."
$ws2.Range("C8").Value = "The entry should be removed from the specified phone book"

$ws2.Range("D3").Value = "Passed"
$ws2.Range("D4").Value = "Passed"
$ws2.Range("A5").Value = "RESTful API endpoint for readAllEntriesFromAllPhoneBooks"
$ws2.Range("B5").Value = "Read the entries from all phone books"
$ws2.Range("D5").Value = "Passed"
$ws2.Range("D6").Value = "Passed"
$ws2.Range("D7").Value = "Passed"
$ws2.Range("D8").Value = "Passed"

# ---------------------------------------------------------------------
# 4. Copy the look & feel (wrap/left/top alignment) of the "Must-Have"
#    sheet onto the used cells of the new sheet.
# ---------------------------------------------------------------------
$ws1.Range("A1").Copy()
$ws2.Range("A1:F1").PasteSpecial(-4122)
$ws2.Range("A2:D8").PasteSpecial(-4122)

# Re-apply the values (PasteSpecial(Formats) does not touch values, but do
# it defensively in case a future engine version clears them)
$ws2.Range("A1").Value = "Requirement"

# ---------------------------------------------------------------------
# 5. Row heights (best-effort match of the autosized wrap-text rows)
# ---------------------------------------------------------------------
$ws2.Rows.Item(1).RowHeight = 51
$ws2.Rows.Item(2).RowHeight = 102
$ws2.Rows.Item(3).RowHeight = 136
$ws2.Rows.Item(4).RowHeight = 136
$ws2.Rows.Item(5).RowHeight = 119
$ws2.Rows.Item(6).RowHeight = 136
$ws2.Rows.Item(7).RowHeight = 119
$ws2.Rows.Item(8).RowHeight = 119

# ---------------------------------------------------------------------
# 6. AutoFilter on the header row + matching hidden defined name
# ---------------------------------------------------------------------
$ws2.Range("A1:F1").AutoFilter()
$fdb = $ws2.Names.Add("_xlnm._FilterDatabase", "=Extra!`$A`$1:`$F`$1")
$fdb.Visible = $false

# ---------------------------------------------------------------------
# 7. View state for the new sheet: freeze header row/col, zoom 140,
#    select A2, and make it the active tab.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("A2").Select()
$excel.ActiveWindow.Zoom = 140

Write-Output "Added 'Extra' worksheet with readAllEntriesFromAllPhoneBooks screen operation."
